$wb = $excel.ActiveWorkbook

# --- Repayment schedule: insert a new (blank) column before the old "Late" column ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of column M so the newly inserted column N inherits it
# (matches the width Excel assigns to a freshly inserted column copying the
# formatting of the column to its left).
$colWidth = $wsSchedule.Columns("M").ColumnWidth

$wsSchedule.Columns("N").Insert()
$wsSchedule.Columns("N").ColumnWidth = $colWidth

# --- Update the active sheet / selections to match the saved view state ---

# "Edit Repayment Schedule" selection moves to B4
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$wsEdit.Activate()
$wsEdit.Range("B4").Select()

# "Repayment schedule" becomes the active tab, with selection on L18
$wsSchedule.Activate()
$wsSchedule.Range("L18").Select()
